$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2666.3333
$ws.Range("I4").Value = 2499.5
$ws.Range("K4").Value = 2499.5
$ws.Range("M4").Value = -2385.5
$ws.Range("H17").Value = 2638.077
$ws.Range("J17").Value = 2144.3333
$ws.Range("L17").Value = 6432.999899999999
$ws.Range("N17").Value = -6768.999899999999
$ws.Range("H39").Value = 66.25
$ws.Range("I39").Value = 66.25
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 198.75
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 97.25
$ws.Range("N39").ClearContents()
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20924
$ws.Range("H86").Value = 3419.8
$ws.Range("I86").Value = 3419.8
$ws.Range("K86").Value = 3419.8
$ws.Range("M86").Value = -2296.8
$ws.Range("H89").Value = 3419.8
$ws.Range("I89").Value = 3419.8
$ws.Range("K89").Value = 17099
$ws.Range("M89").Value = -11483
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4446.4067
$ws.Range("I32").Value = 3467.5881
$ws.Range("K32").Value = 3467.5881
$ws.Range("M32").Value = -3180.5881
$ws.Range("H37").Value = 20000
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20546
$ws.Range("H74").Value = 1012.0732
$ws.Range("I74").Value = 520.7353000000001
$ws.Range("K74").Value = 520.7353000000001
$ws.Range("M74").Value = 353.2646999999999
$ws.Range("H77").Value = 1012.0732
$ws.Range("I77").Value = 520.7353000000001
$ws.Range("K77").Value = 2603.6765
$ws.Range("M77").Value = 1764.3235
$ws.Range("H97").Value = 1021.6842
$ws.Range("I97").Value = 944.64703
$ws.Range("K97").Value = 944.64703
$ws.Range("M97").Value = -448.64703
$ws.Range("H113").Value = 50000
$ws.Range("J113").Value = 50000
$ws.Range("L113").Value = 50000
$ws.Range("N113").Value = -58678
$ws.Range("H132").Value = 1378.9524
$ws.Range("I132").Value = 1339.2106
$ws.Range("J132").Value = 1756.5
$ws.Range("K132").Value = 4017.6318
$ws.Range("L132").Value = 5269.5
$ws.Range("M132").Value = -1487.6318
$ws.Range("N132").Value = -10329.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 823.55554
$ws.Range("I94").Value = 823.55554
$ws.Range("K94").Value = 823.55554
$ws.Range("M94").Value = -372.55554
$ws.Range("H105").Value = 2543.0908
$ws.Range("I105").Value = 2543.0908
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2543.0908
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -796.0907999999999
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 5845.2144
$ws.Range("I134").Value = 6455.7085
$ws.Range("J134").Value = 2182.25
$ws.Range("K134").Value = 19367.1255
$ws.Range("L134").Value = 6546.75
$ws.Range("M134").Value = -16832.1255
$ws.Range("N134").Value = -11616.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 86446.28999999999
$ws.Range("I16").Value = 100751.836
$ws.Range("J16").Value = 613
$ws.Range("K16").Value = 100751.836
$ws.Range("L16").Value = 613
$ws.Range("M16").Value = -100464.836
$ws.Range("N16").Value = -1187
$ws.Range("H52").Value = 65280
$ws.Range("J52").Value = 65280
$ws.Range("L52").Value = 65280
$ws.Range("N52").Value = -65868
$ws.Range("H58").Value = 4351384.5
$ws.Range("I58").Value = 14493283
$ws.Range("K58").Value = 14493283
$ws.Range("M58").Value = -14493080
$ws.Range("H113").Value = 86446.28999999999
$ws.Range("I113").Value = 100751.836
$ws.Range("J113").Value = 613
$ws.Range("K113").Value = 100751.836
$ws.Range("L113").Value = 613
$ws.Range("M113").Value = -98581.836
$ws.Range("N113").Value = -4953
$ws.Range("H122").Value = 5638.625
$ws.Range("I122").Value = 4349.5
$ws.Range("K122").Value = 13048.5
$ws.Range("M122").Value = -10598.5
$ws.Range("H132").Value = 5805
$ws.Range("I132").Value = 4999
$ws.Range("J132").Value = 6342.3335
$ws.Range("K132").Value = 14997
$ws.Range("L132").Value = 19027.0005
$ws.Range("M132").Value = -12467
$ws.Range("N132").Value = -24087.0005
$ws.Range("H136").Value = 4351384.5
$ws.Range("I136").Value = 14493283
$ws.Range("K136").Value = 43479849
$ws.Range("M136").Value = -43477299

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 481.9375
$ws.Range("J107").Value = 485.46155
$ws.Range("L107").Value = 1456.38465
$ws.Range("N107").Value = -5296.38465
$ws.Range("H113").Value = 7966.0713
$ws.Range("J113").Value = 820.36365
$ws.Range("L113").Value = 2461.09095
$ws.Range("N113").Value = -6801.09095
$ws.Range("H114").Value = 1094.7778
$ws.Range("I114").Value = 307
$ws.Range("J114").Value = 1319.8572
$ws.Range("K114").Value = 921
$ws.Range("L114").Value = 3959.5716
$ws.Range("M114").Value = 2333
$ws.Range("N114").Value = -10467.5716
$ws.Range("H131").Value = 16602.137
$ws.Range("I131").Value = 760
$ws.Range("J131").Value = 17761.316
$ws.Range("K131").Value = 2280
$ws.Range("L131").Value = 53283.948
$ws.Range("M131").Value = 2760
$ws.Range("N131").Value = -63363.948
$ws.Range("H132").Value = 1974.25
$ws.Range("I132").Value = 1001
$ws.Range("J132").Value = 2113.2856
$ws.Range("K132").Value = 9009
$ws.Range("L132").Value = 19019.5704
$ws.Range("M132").Value = -6479
$ws.Range("N132").Value = -24079.5704
$ws.Range("H140").Value = 2269.8928
$ws.Range("I140").Value = 1425.2
$ws.Range("J140").Value = 3244.5386
$ws.Range("K140").Value = 4275.6
$ws.Range("L140").Value = 9733.6158
$ws.Range("M140").Value = 904.3999999999996
$ws.Range("N140").Value = -20093.6158

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 34758.2
$ws.Range("J15").Value = 34758.2
$ws.Range("L15").Value = 34758.2
$ws.Range("N15").Value = -35334.2
$ws.Range("H70").Value = 4271.8
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4271.8
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4271.8
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -4811.8
$ws.Range("H73").Value = 4271.8
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4271.8
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4271.8
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6143.8
$ws.Range("H81").Value = 34758.2
$ws.Range("J81").Value = 34758.2
$ws.Range("L81").Value = 34758.2
$ws.Range("N81").Value = -36754.2
$ws.Range("H84").Value = 34758.2
$ws.Range("J84").Value = 34758.2
$ws.Range("L84").Value = 104274.6
$ws.Range("N84").Value = -114258.6
$ws.Range("H97").Value = 2634.2856
$ws.Range("I97").Value = 2573.3333
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 2573.3333
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -2077.3333
$ws.Range("N97").Value = -3992
$ws.Range("H126").Value = 1573628.4
$ws.Range("I126").Value = 4632612.5
$ws.Range("J126").Value = 44136.375
$ws.Range("K126").Value = 13897837.5
$ws.Range("L126").Value = 132409.125
$ws.Range("M126").Value = -13895367.5
$ws.Range("N126").Value = -137349.125

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2773.6
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5404
$ws.Range("H113").Value = 2773.6
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 9749.875
$ws.Range("I122").Value = 10799.8
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 32399.4
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -29949.4
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 4306.3335
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4306.3335
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12919.0005
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -17979.0005
$ws.Range("H136").Value = 5463.1177
$ws.Range("I136").Value = 4179.7
$ws.Range("J136").Value = 7296.5713
$ws.Range("K136").Value = 12539.1
$ws.Range("L136").Value = 21889.7139
$ws.Range("M136").Value = -9989.099999999999
$ws.Range("N136").Value = -26989.7139

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 768.6
$ws.Range("I100").Value = 585.75
$ws.Range("K100").Value = 1171.5
$ws.Range("M100").Value = -630.5
$ws.Range("H108").Value = 24500
$ws.Range("J108").Value = 24500
$ws.Range("L108").Value = 24500
$ws.Range("N108").Value = -32180
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22954
$ws.Range("H126").Value = 12930.5
$ws.Range("I126").Value = 14763.25
$ws.Range("J126").Value = 5599.5
$ws.Range("K126").Value = 44289.75
$ws.Range("L126").Value = 16798.5
$ws.Range("M126").Value = -41819.75
$ws.Range("N126").Value = -21738.5
$ws.Range("H132").Value = 1384.303
$ws.Range("I132").Value = 926.2593000000001
$ws.Range("J132").Value = 3445.5
$ws.Range("K132").Value = 2778.7779
$ws.Range("L132").Value = 10336.5
$ws.Range("M132").Value = -248.7779
$ws.Range("N132").Value = -15396.5
$ws.Range("H136").Value = 16836842
$ws.Range("I136").Value = 30866326
$ws.Range("J136").Value = 1459.6
$ws.Range("K136").Value = 92598978
$ws.Range("L136").Value = 4378.799999999999
$ws.Range("M136").Value = -92596428
$ws.Range("N136").Value = -9478.799999999999
